# Updated Phenology DAS in observed files
#
# The Cotton.Phenology.*DAS columns (K:P = Emergence/Squaring/Flowering/
# Cutout/Openbolls/Maturity/HarvestRipeDAS) only had a single column filled
# in per growth-stage row. Once a stage has been reached, every later sample
# within that sowing's block should report the DAS of each stage that has
# already occurred by that point (K..P line up "as of this sample").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Emerald2015P1 block (rows 2-8) ---
$ws.Range("K2:K8").Value = 65
$ws.Range("L2:L8").Value = 85
$ws.Range("M2:M8").Value = 121
$ws.Range("N2:N8").Value = 136
$ws.Range("O2:O8").Value = 159
$ws.Range("P2:P8").Value = 177

# --- Emerald2015P2 block (rows 9-15) ---
$ws.Range("K9:K15").Value = 56
$ws.Range("L9:L15").Value = 71
$ws.Range("M9:M15").Value = 106
$ws.Range("N9:N15").Value = 122
$ws.Range("O9:O15").Value = 143
$ws.Range("P9:P15").Value = 160

# --- Emerald2015P3 block (rows 16-22) ---
$ws.Range("K16:K22").Value = 49
$ws.Range("L16:L22").Value = 71
$ws.Range("M16:M22").Value = 107
$ws.Range("N16:N22").Value = 124
$ws.Range("O16:O22").Value = 140
$ws.Range("P16:P22").Value = 161

# Reflect the author's final on-screen selection: the last thing touched
# was the bottom-right (scrolling) pane, ending with P16:P22 selected.
$ws.Range("P16:P22").Select()
